$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-23 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-24 Monday", 2) | Out-Null
$d.Content.Find.Execute("76-73=", $true, $false, $false, $false, $false, $true, 1, $false, "16+81=", 2) | Out-Null
$d.Content.Find.Execute("56-54=", $true, $false, $false, $false, $false, $true, 1, $false, "49-11=", 2) | Out-Null
$d.Content.Find.Execute("30+38=", $true, $false, $false, $false, $false, $true, 1, $false, "60-54=", 2) | Out-Null
$d.Content.Find.Execute("57-12=", $true, $false, $false, $false, $false, $true, 1, $false, "21+28=", 2) | Out-Null
$d.Content.Find.Execute("51+14=", $true, $false, $false, $false, $false, $true, 1, $false, "35+32=", 2) | Out-Null
$d.Content.Find.Execute("3+23=", $true, $false, $false, $false, $false, $true, 1, $false, "12+55=", 2) | Out-Null
$d.Content.Find.Execute("38+11=", $true, $false, $false, $false, $false, $true, 1, $false, "8+61=", 2) | Out-Null
$d.Content.Find.Execute("9+3=", $true, $false, $false, $false, $false, $true, 1, $false, "18+41=", 2) | Out-Null
$d.Content.Find.Execute("87-16=", $true, $false, $false, $false, $false, $true, 1, $false, "7+12=", 2) | Out-Null
$d.Content.Find.Execute("65+29=", $true, $false, $false, $false, $false, $true, 1, $false, "66+5=", 2) | Out-Null
$d.Content.Find.Execute("95-34=", $true, $false, $false, $false, $false, $true, 1, $false, "89-13=", 2) | Out-Null
$d.Content.Find.Execute("16+15=", $true, $false, $false, $false, $false, $true, 1, $false, "52-23=", 2) | Out-Null
$d.Content.Find.Execute("93-84=", $true, $false, $false, $false, $false, $true, 1, $false, "99-48=", 2) | Out-Null
$d.Content.Find.Execute("62-58=", $true, $false, $false, $false, $false, $true, 1, $false, "53-0=", 2) | Out-Null
$d.Content.Find.Execute("5+47=", $true, $false, $false, $false, $false, $true, 1, $false, "35+29=", 2) | Out-Null
$d.Content.Find.Execute("34-17=", $true, $false, $false, $false, $false, $true, 1, $false, "37-3=", 2) | Out-Null
$d.Content.Find.Execute("53-40=", $true, $false, $false, $false, $false, $true, 1, $false, "17+82=", 2) | Out-Null
$d.Content.Find.Execute("34+16=", $true, $false, $false, $false, $false, $true, 1, $false, "19+45=", 2) | Out-Null
$d.Content.Find.Execute("78-31=", $true, $false, $false, $false, $false, $true, 1, $false, "84+11=", 2) | Out-Null
$d.Content.Find.Execute("4+87=", $true, $false, $false, $false, $false, $true, 1, $false, "64+16=", 2) | Out-Null
$d.Content.Find.Execute("60-16=", $true, $false, $false, $false, $false, $true, 1, $false, "25+53=", 2) | Out-Null
$d.Content.Find.Execute("39+13=", $true, $false, $false, $false, $false, $true, 1, $false, "65-50=", 2) | Out-Null
$d.Content.Find.Execute("49-21=", $true, $false, $false, $false, $false, $true, 1, $false, "65-33=", 2) | Out-Null
$d.Content.Find.Execute("21+26=", $true, $false, $false, $false, $false, $true, 1, $false, "28-24=", 2) | Out-Null
$d.Content.Find.Execute("67-0=", $true, $false, $false, $false, $false, $true, 1, $false, "93-71=", 2) | Out-Null
$d.Content.Find.Execute("58-56=", $true, $false, $false, $false, $false, $true, 1, $false, "35+51=", 2) | Out-Null
$d.Content.Find.Execute("73-11=", $true, $false, $false, $false, $false, $true, 1, $false, "15+47=", 2) | Out-Null
$d.Content.Find.Execute("76-5=", $true, $false, $false, $false, $false, $true, 1, $false, "75-36=", 2) | Out-Null
$d.Content.Find.Execute("81-80=", $true, $false, $false, $false, $false, $true, 1, $false, "15-3=", 2) | Out-Null
$d.Content.Find.Execute("97-4=", $true, $false, $false, $false, $false, $true, 1, $false, "64-57=", 2) | Out-Null
$d.Content.Find.Execute("89-77=", $true, $false, $false, $false, $false, $true, 1, $false, "69+13=", 2) | Out-Null
$d.Content.Find.Execute("47-20=", $true, $false, $false, $false, $false, $true, 1, $false, "44+26=", 2) | Out-Null
$d.Content.Find.Execute("61-30=", $true, $false, $false, $false, $false, $true, 1, $false, "54-39=", 2) | Out-Null
$d.Content.Find.Execute("25-11=", $true, $false, $false, $false, $false, $true, 1, $false, "46+35=", 2) | Out-Null
$d.Content.Find.Execute("61-24=", $true, $false, $false, $false, $false, $true, 1, $false, "14+62=", 2) | Out-Null
$d.Content.Find.Execute("60-12=", $true, $false, $false, $false, $false, $true, 1, $false, "85-41=", 2) | Out-Null
$d.Content.Find.Execute("70+4=", $true, $false, $false, $false, $false, $true, 1, $false, "22+47=", 2) | Out-Null
$d.Content.Find.Execute("76+18=", $true, $false, $false, $false, $false, $true, 1, $false, "42+54=", 2) | Out-Null
$d.Content.Find.Execute("56+14=", $true, $false, $false, $false, $false, $true, 1, $false, "16+75=", 2) | Out-Null
$d.Content.Find.Execute("68-22=", $true, $false, $false, $false, $false, $true, 1, $false, "77+2=", 2) | Out-Null
$d.Content.Find.Execute("8+55=", $true, $false, $false, $false, $false, $true, 1, $false, "46+29=", 2) | Out-Null
$d.Content.Find.Execute("77-15=", $true, $false, $false, $false, $false, $true, 1, $false, "46-39=", 2) | Out-Null
$d.Content.Find.Execute("60+36=", $true, $false, $false, $false, $false, $true, 1, $false, "2+44=", 2) | Out-Null
$d.Content.Find.Execute("61-13=", $true, $false, $false, $false, $false, $true, 1, $false, "77-75=", 2) | Out-Null
$d.Content.Find.Execute("82-74=", $true, $false, $false, $false, $false, $true, 1, $false, "90-48=", 2) | Out-Null
$d.Content.Find.Execute("10+8=", $true, $false, $false, $false, $false, $true, 1, $false, "42+21=", 2) | Out-Null
$d.Content.Find.Execute("31-12=", $true, $false, $false, $false, $false, $true, 1, $false, "60-10=", 2) | Out-Null
$d.Content.Find.Execute("66-43=", $true, $false, $false, $false, $false, $true, 1, $false, "28+59=", 2) | Out-Null
$d.Content.Find.Execute("63-30=", $true, $false, $false, $false, $false, $true, 1, $false, "49-11=", 2) | Out-Null
$d.Content.Find.Execute("14+67=", $true, $false, $false, $false, $false, $true, 1, $false, "11+60=", 2) | Out-Null
$d.Content.Find.Execute("49-30=", $true, $false, $false, $false, $false, $true, 1, $false, "6+0=", 2) | Out-Null
$d.Content.Find.Execute("28+66=", $true, $false, $false, $false, $false, $true, 1, $false, "81-64=", 2) | Out-Null
$d.Content.Find.Execute("51-16=", $true, $false, $false, $false, $false, $true, 1, $false, "89-29=", 2) | Out-Null
$d.Content.Find.Execute("97-95=", $true, $false, $false, $false, $false, $true, 1, $false, "35-33=", 2) | Out-Null
$d.Content.Find.Execute("89-71=", $true, $false, $false, $false, $false, $true, 1, $false, "72-51=", 2) | Out-Null
$d.Content.Find.Execute("64-61=", $true, $false, $false, $false, $false, $true, 1, $false, "57-26=", 2) | Out-Null
$d.Content.Find.Execute("7+74=", $true, $false, $false, $false, $false, $true, 1, $false, "12+41=", 2) | Out-Null
$d.Content.Find.Execute("90-6=", $true, $false, $false, $false, $false, $true, 1, $false, "56-53=", 2) | Out-Null
$d.Content.Find.Execute("76-9=", $true, $false, $false, $false, $false, $true, 1, $false, "69-14=", 2) | Out-Null
$d.Content.Find.Execute("76+20=", $true, $false, $false, $false, $false, $true, 1, $false, "21-9=", 2) | Out-Null
$d.Content.Find.Execute("10-8=", $true, $false, $false, $false, $false, $true, 1, $false, "80-41=", 2) | Out-Null
$d.Content.Find.Execute("35-31=", $true, $false, $false, $false, $false, $true, 1, $false, "93-82=", 2) | Out-Null
$d.Content.Find.Execute("17+8=", $true, $false, $false, $false, $false, $true, 1, $false, "72-5=", 2) | Out-Null
$d.Content.Find.Execute("48-40=", $true, $false, $false, $false, $false, $true, 1, $false, "66+21=", 2) | Out-Null
$d.Content.Find.Execute("87-83=", $true, $false, $false, $false, $false, $true, 1, $false, "32-23=", 2) | Out-Null
$d.Content.Find.Execute("6+62=", $true, $false, $false, $false, $false, $true, 1, $false, "81-25=", 2) | Out-Null
$d.Content.Find.Execute("89-37=", $true, $false, $false, $false, $false, $true, 1, $false, "9+8=", 2) | Out-Null
$d.Content.Find.Execute("3+19=", $true, $false, $false, $false, $false, $true, 1, $false, "99-25=", 2) | Out-Null
$d.Content.Find.Execute("71-10=", $true, $false, $false, $false, $false, $true, 1, $false, "70-18=", 2) | Out-Null
$d.Content.Find.Execute("53-15=", $true, $false, $false, $false, $false, $true, 1, $false, "88+2=", 2) | Out-Null
$d.Content.Find.Execute("0-0=", $true, $false, $false, $false, $false, $true, 1, $false, "11+28=", 2) | Out-Null
$d.Content.Find.Execute("45+17=", $true, $false, $false, $false, $false, $true, 1, $false, "98-61=", 2) | Out-Null
$d.Content.Find.Execute("80+18=", $true, $false, $false, $false, $false, $true, 1, $false, "77+10=", 2) | Out-Null
$d.Content.Find.Execute("0+10=", $true, $false, $false, $false, $false, $true, 1, $false, "19+22=", 2) | Out-Null
$d.Content.Find.Execute("96-69=", $true, $false, $false, $false, $false, $true, 1, $false, "90-88=", 2) | Out-Null
$d.Content.Find.Execute("76-7=", $true, $false, $false, $false, $false, $true, 1, $false, "22+4=", 2) | Out-Null
$d.Content.Find.Execute("19+75=", $true, $false, $false, $false, $false, $true, 1, $false, "65-60=", 2) | Out-Null
$d.Content.Find.Execute("10+62=", $true, $false, $false, $false, $false, $true, 1, $false, "85-66=", 2) | Out-Null
$d.Content.Find.Execute("66-15=", $true, $false, $false, $false, $false, $true, 1, $false, "71-41=", 2) | Out-Null
$d.Content.Find.Execute("98-45=", $true, $false, $false, $false, $false, $true, 1, $false, "74-4=", 2) | Out-Null
$d.Content.Find.Execute("72+24=", $true, $false, $false, $false, $false, $true, 1, $false, "3-2=", 2) | Out-Null
$d.Content.Find.Execute("38+1=", $true, $false, $false, $false, $false, $true, 1, $false, "6+84=", 2) | Out-Null
$d.Content.Find.Execute("48+19=", $true, $false, $false, $false, $false, $true, 1, $false, "63+26=", 2) | Out-Null
$d.Content.Find.Execute("43-31=", $true, $false, $false, $false, $false, $true, 1, $false, "28+11=", 2) | Out-Null
$d.Content.Find.Execute("44+38=", $true, $false, $false, $false, $false, $true, 1, $false, "27-19=", 2) | Out-Null
$d.Content.Find.Execute("83+16=", $true, $false, $false, $false, $false, $true, 1, $false, "12+5=", 2) | Out-Null
$d.Content.Find.Execute("60-20=", $true, $false, $false, $false, $false, $true, 1, $false, "27-11=", 2) | Out-Null
$d.Content.Find.Execute("20+21=", $true, $false, $false, $false, $false, $true, 1, $false, "2+77=", 2) | Out-Null
$d.Content.Find.Execute("13+51=", $true, $false, $false, $false, $false, $true, 1, $false, "96-40=", 2) | Out-Null
$d.Content.Find.Execute("9+7=", $true, $false, $false, $false, $false, $true, 1, $false, "37+33=", 2) | Out-Null
$d.Content.Find.Execute("86-4=", $true, $false, $false, $false, $false, $true, 1, $false, "40-32=", 2) | Out-Null
$d.Content.Find.Execute("73-37=", $true, $false, $false, $false, $false, $true, 1, $false, "34+11=", 2) | Out-Null
$d.Content.Find.Execute("0+79=", $true, $false, $false, $false, $false, $true, 1, $false, "77-38=", 2) | Out-Null
$d.Content.Find.Execute("96-27=", $true, $false, $false, $false, $false, $true, 1, $false, "52-19=", 2) | Out-Null
$d.Content.Find.Execute("60-47=", $true, $false, $false, $false, $false, $true, 1, $false, "73-9=", 2) | Out-Null
$d.Content.Find.Execute("38+61=", $true, $false, $false, $false, $false, $true, 1, $false, "9+16=", 2) | Out-Null
$d.Content.Find.Execute("95-31=", $true, $false, $false, $false, $false, $true, 1, $false, "4+10=", 2) | Out-Null
$d.Content.Find.Execute("44-32=", $true, $false, $false, $false, $false, $true, 1, $false, "46-23=", 2) | Out-Null
$d.Content.Find.Execute("75-18=", $true, $false, $false, $false, $false, $true, 1, $false, "95-94=", 2) | Out-Null
$d.Content.Find.Execute("7+75=", $true, $false, $false, $false, $false, $true, 1, $false, "47-42=", 2) | Out-Null
